$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

# B11 currently holds the text "R40" (rule name); it needs to become the
# text "1" while keeping its existing style/number-format untouched.
#
# A plain `Range("B11").Value = "1"` would make Excel auto-detect the
# numeric-looking literal and store it as a number (and forcing text via
# NumberFormat="@" on B11 itself would stamp a brand-new text number
# format onto its style, which the target file does not have).
#
# Instead: stage the literal "1" as real text on a scratch cell far away
# (format it as Text so Excel keeps the leading content as a string),
# copy just the VALUE (not formatting) onto B11 so its original style id
# is preserved, then remove the scratch row entirely so nothing else in
# the workbook is touched.
$scratch = $ws.Range("A200")
$scratch.NumberFormat = "@"
$scratch.Value = "1"
$scratch.Copy()

$ws.Range("B11").PasteSpecial(-4163) # xlPasteValues

$scratch.EntireRow.Delete()

$excel.CutCopyMode = $false
$wb.Save()
